$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# The "ultimate" row (row 5: a zero-paid placeholder row dated 40912, tagged
# with the now-obsolete "ultimate" label) is removed entirely, shifting the
# remaining data rows (old rows 6-10, the J49:L52 projection rows) up by one.
$ws.Rows.Item(5).Delete() | Out-Null

# Leave the selection on the row that now occupies the old row 5's position,
# matching the post-edit UI state (whole row highlighted).
$ws.Range("A5:XFD5").Select() | Out-Null
